$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "63.839.37"
Set-TextValue $ws.Range("E2") "  +1.25%  "

Set-TextValue $ws.Range("D3") "3.307.68"
Set-TextValue $ws.Range("E3") "  +5.83%  "

Set-TextValue $ws.Range("E4") "  +0.03%  "

Set-TextValue $ws.Range("D5") "602.70"
Set-TextValue $ws.Range("E5") "  +1.08%  "

Set-TextValue $ws.Range("D6") "142.65"
Set-TextValue $ws.Range("E6") "  +4.43%  "

Set-TextValue $ws.Range("E7") "  -0.01%  "

Set-TextValue $ws.Range("D8") "3.305.06"
Set-TextValue $ws.Range("E8") "  +5.98%  "

Set-TextValue $ws.Range("E9") "  +1.16%  "

Set-TextValue $ws.Range("E10") "  +2.46%  "

Set-TextValue $ws.Range("D11") "5.52"
Set-TextValue $ws.Range("E11") "  +4.52%  "

Set-TextValue $ws.Range("D12") "0.472"
Set-TextValue $ws.Range("E12") "  +3.43%  "

Set-TextValue $ws.Range("D13") "0.0000248"
Set-TextValue $ws.Range("E13") "  +0.66%  "

Set-TextValue $ws.Range("D14") "34.74"
Set-TextValue $ws.Range("E14") "  +1.25%  "

Set-TextValue $ws.Range("D15") "3.855.43"
Set-TextValue $ws.Range("E15") "  +6.01%  "

Set-TextValue $ws.Range("E16") "  +0.25%  "

Set-TextValue $ws.Range("D17") "3.309.47"
Set-TextValue $ws.Range("E17") "  +6.10%  "

Set-TextValue $ws.Range("D18") "63.923.80"
Set-TextValue $ws.Range("E18") "  +1.36%  "

Set-TextValue $ws.Range("D19") "6.89"
Set-TextValue $ws.Range("E19") "  +2.99%  "

Set-TextValue $ws.Range("D20") "481.03"
Set-TextValue $ws.Range("E20") "  +1.49%  "

Set-TextValue $ws.Range("D21") "14.20"
Set-TextValue $ws.Range("E21") "  +0.34%  "

Set-TextValue $ws.Range("D22") "0.735"
Set-TextValue $ws.Range("E22") "  +5.26%  "

Set-TextValue $ws.Range("D23") "8.01"
Set-TextValue $ws.Range("E23") "  +4.26%  "

Set-TextValue $ws.Range("D24") "85.06"
Set-TextValue $ws.Range("E24") "  -1.28%  "

Set-TextValue $ws.Range("D25") "13.46"
Set-TextValue $ws.Range("E25") "  +3.74%  "

Set-TextValue $ws.Range("E26") "  -0.02%  "

Set-TextValue $ws.Range("E27") "  +1.69%  "

Set-TextValue $ws.Range("D28") "7.28"
Set-TextValue $ws.Range("E28") "  +4.40%  "

Set-TextValue $ws.Range("E29") "  -0.01%  "

Set-TextValue $ws.Range("D30") "8.15"
Set-TextValue $ws.Range("E30") "  +2.74%  "

Set-TextValue $ws.Range("E31") "  +3.76%  "

Set-TextValue $ws.Range("D32") "29.40"
Set-TextValue $ws.Range("E32") "  +9.90%  "

Set-TextValue $ws.Range("E33") "  -2.10%  "

Set-TextValue $ws.Range("E34") "  +0.46%  "

Set-TextValue $ws.Range("E35") "  +1.78%  "

Set-TextValue $ws.Range("D36") "5.97"
Set-TextValue $ws.Range("E36") "  +2.75%  "

Set-TextValue $ws.Range("D37") "52.82"
Set-TextValue $ws.Range("E37") "  +1.61%  "

Set-TextValue $ws.Range("D38") "0.0₃0750"
Set-TextValue $ws.Range("E38") "  +6.71%  "

Set-TextValue $ws.Range("D39") "0.0403"
Set-TextValue $ws.Range("E39") "  +3.94%  "

Set-TextValue $ws.Range("D40") "429.73"
Set-TextValue $ws.Range("E40") "  +2.05%  "

Set-TextValue $ws.Range("D41") "3.041.96"
Set-TextValue $ws.Range("E41") "  +4.87%  "

Set-TextValue $ws.Range("D42") "8.40"
Set-TextValue $ws.Range("E42") "  +2.07%  "

Set-TextValue $ws.Range("E43") "  +1.92%  "

Set-TextValue $ws.Range("E44") "  -0.88%  "

Set-TextValue $ws.Range("E45") "  -0.60%  "

Set-TextValue $ws.Range("D46") "2.20"
Set-TextValue $ws.Range("E46") "  +3.56%  "

Set-TextValue $ws.Range("D47") "26.39"
Set-TextValue $ws.Range("E47") "  +3.13%  "

Set-TextValue $ws.Range("E48") "  +0.03%  "

Set-TextValue $ws.Range("D49") "35.62"
Set-TextValue $ws.Range("E49") "  +11.04%  "

Set-TextValue $ws.Range("E50") "  +1.98%  "

Set-TextValue $ws.Range("D51") "2.31"
Set-TextValue $ws.Range("E51") "  +2.44%  "
